# Refresh cached market-board figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# for the rows whose quoted prices moved in this run. Values below are the new
# authoritative figures pulled by the scheduled pricing job.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21 (Leve Item ID 2149)
$ws.Range("H21").Value = 27839
$ws.Range("I21").Value = 1758.5
$ws.Range("J21").Value = 80000
$ws.Range("K21").Value = 1758.5
$ws.Range("L21").Value = 80000
$ws.Range("M21").Value = -1290.5
$ws.Range("N21").Value = -80936
# Row 23 (Leve Item ID 2149)
$ws.Range("H23").Value = 27839
$ws.Range("I23").Value = 1758.5
$ws.Range("J23").Value = 80000
$ws.Range("K23").Value = 1758.5
$ws.Range("L23").Value = 80000
$ws.Range("M23").Value = -1524.5
$ws.Range("N23").Value = -80468
# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 337.34784
$ws.Range("I33").Value = 339.04544
$ws.Range("K33").Value = 339.04544
$ws.Range("M33").Value = -110.04544
# Row 34 (Leve Item ID 2160)
$ws.Range("H34").Value = 870.5714
$ws.Range("I34").Value = 870.5714
$ws.Range("K34").Value = 870.5714
$ws.Range("M34").Value = -667.5714
# Row 36 (Leve Item ID 2160)
$ws.Range("H36").Value = 870.5714
$ws.Range("I36").Value = 870.5714
$ws.Range("K36").Value = 870.5714
$ws.Range("M36").Value = -155.5714
# Row 47 (Leve Item ID 2169)
$ws.Range("H47").Value = 13018.5
$ws.Range("I47").Value = 11666.667
$ws.Range("J47").Value = 17074
$ws.Range("K47").Value = 11666.667
$ws.Range("L47").Value = 17074
$ws.Range("M47").Value = -10694.667
$ws.Range("N47").Value = -19018
# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 2962041.5
$ws.Range("I64").Value = 5131482
$ws.Range("J64").Value = 3713.6365
$ws.Range("K64").Value = 5131482
$ws.Range("L64").Value = 3713.6365
$ws.Range("M64").Value = -5131234
$ws.Range("N64").Value = -4209.636500000001
# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 2962041.5
$ws.Range("I67").Value = 5131482
$ws.Range("J67").Value = 3713.6365
$ws.Range("K67").Value = 5131482
$ws.Range("L67").Value = 3713.6365
$ws.Range("M67").Value = -5130624
$ws.Range("N67").Value = -5429.636500000001
# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 1207.96
$ws.Range("I129").Value = 388.17648
$ws.Range("K129").Value = 1164.52944
$ws.Range("M129").Value = 3835.47056
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 4651850
$ws.Range("I137").Value = 619.5714
$ws.Range("J137").Value = 13334146
$ws.Range("K137").Value = 1858.7142
$ws.Range("L137").Value = 40002438
$ws.Range("M137").Value = 691.2857999999999
$ws.Range("N137").Value = -40007538

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 7894.228
$ws.Range("I32").Value = 8099.644
$ws.Range("J32").Value = 7105
$ws.Range("K32").Value = 8099.644
$ws.Range("L32").Value = 7105
$ws.Range("M32").Value = -7812.644
$ws.Range("N32").Value = -7679
# Row 44 (Leve Item ID 3861)
$ws.Range("H44").Value = 29900
$ws.Range("J44").Value = 29900
$ws.Range("L44").Value = 29900
$ws.Range("N44").Value = -30876
# Row 51 (Leve Item ID 3858)
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 10205472
$ws.Range("I61").Value = 11629378
$ws.Range("J61").Value = 810.3333
$ws.Range("K61").Value = 11629378
$ws.Range("L61").Value = 810.3333
$ws.Range("M61").Value = -11629166
$ws.Range("N61").Value = -1234.3333
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 6946665.5
$ws.Range("I132").Value = 8930575
$ws.Range("J132").Value = 2982
$ws.Range("K132").Value = 26791725
$ws.Range("L132").Value = 8946
$ws.Range("M132").Value = -26789195
$ws.Range("N132").Value = -14006
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 10205472
$ws.Range("I136").Value = 11629378
$ws.Range("J136").Value = 810.3333
$ws.Range("K136").Value = 34888134
$ws.Range("L136").Value = 2430.9999
$ws.Range("M136").Value = -34885584
$ws.Range("N136").Value = -7530.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1909.3562
$ws.Range("I134").Value = 1174.8448
$ws.Range("K134").Value = 3524.5344
$ws.Range("M134").Value = -989.5344000000005

$ws = $wb.Worksheets.Item("CRP")
# Row 26 (Leve Item ID 2004)
$ws.Range("H26").Value = 10010.272
$ws.Range("I26").Value = 1256.25
$ws.Range("J26").Value = 33354.332
$ws.Range("K26").Value = 1256.25
$ws.Range("L26").Value = 33354.332
$ws.Range("M26").Value = -969.25
$ws.Range("N26").Value = -33928.332
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 6668667
$ws.Range("I31").Value = 1808.8937
$ws.Range("K31").Value = 1808.8937
$ws.Range("M31").Value = -1513.8937
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 6668667
$ws.Range("I34").Value = 1808.8937
$ws.Range("K34").Value = 1808.8937
$ws.Range("M34").Value = -1606.8937
# Row 44 (Leve Item ID 1850)
$ws.Range("H44").Value = 114000
$ws.Range("I44").Value = 114000
$ws.Range("K44").Value = 114000
$ws.Range("M44").Value = -113558
# Row 54 (Leve Item ID 2413)
$ws.Range("H54").Value = 30069
$ws.Range("J54").Value = 30069
$ws.Range("L54").Value = 30069
$ws.Range("N54").Value = -31385
# Row 59 (Leve Item ID 1942)
$ws.Range("H59").Value = 13649.667
$ws.Range("I59").Value = 4000
$ws.Range("J59").Value = 15579.6
$ws.Range("K59").Value = 4000
$ws.Range("L59").Value = 15579.6
$ws.Range("M59").Value = -2855
$ws.Range("N59").Value = -17869.6
# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 2239.6365
$ws.Range("I62").Value = 2197.7778
$ws.Range("K62").Value = 2197.7778
$ws.Range("M62").Value = -1573.7778
# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 2239.6365
$ws.Range("I65").Value = 2197.7778
$ws.Range("K65").Value = 10988.889
$ws.Range("M65").Value = -7868.888999999999
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 2176.8518
$ws.Range("I122").Value = 2257.5
$ws.Range("J122").Value = 1946.4286
$ws.Range("K122").Value = 6772.5
$ws.Range("L122").Value = 5839.2858
$ws.Range("M122").Value = -4322.5
$ws.Range("N122").Value = -10739.2858
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 1480.2222
$ws.Range("I134").Value = 1375.74
$ws.Range("K134").Value = 4127.22
$ws.Range("M134").Value = -1592.22

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 778.8095
$ws.Range("I5").Value = 345.33334
$ws.Range("J5").Value = 1356.7778
$ws.Range("K5").Value = 1036.00002
$ws.Range("L5").Value = 4070.3334
$ws.Range("M5").Value = -924.0000199999999
$ws.Range("N5").Value = -4294.3334
# Row 12 (Leve Item ID 4854)
$ws.Range("H12").Value = 29.190475
$ws.Range("I12").Value = 14.7
$ws.Range("J12").Value = 42.363636
$ws.Range("K12").Value = 44.09999999999999
$ws.Range("L12").Value = 127.090908
$ws.Range("M12").Value = 128.9
$ws.Range("N12").Value = -473.090908
# Row 42 (Leve Item ID 4670)
$ws.Range("H42").Value = 2465.5
$ws.Range("I42").Value = 600
$ws.Range("J42").Value = 3398.25
$ws.Range("K42").Value = 1800
$ws.Range("L42").Value = 10194.75
$ws.Range("M42").Value = -1266
$ws.Range("N42").Value = -11262.75
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 2326.5
$ws.Range("J122").Value = 1402.5
$ws.Range("L122").Value = 12622.5
$ws.Range("N122").Value = -17522.5
# Row 134 (Leve Item ID 44074)
$ws.Range("H134").Value = 3415.9707
$ws.Range("I134").Value = 1814.091
$ws.Range("J134").Value = 6352.75
$ws.Range("K134").Value = 5442.272999999999
$ws.Range("L134").Value = 19058.25
$ws.Range("M134").Value = -372.2729999999992
$ws.Range("N134").Value = -29198.25
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 778.8095
$ws.Range("I135").Value = 345.33334
$ws.Range("J135").Value = 1356.7778
$ws.Range("K135").Value = 3108.00006
$ws.Range("L135").Value = 12211.0002
$ws.Range("M135").Value = -573.0000600000003
$ws.Range("N135").Value = -17281.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 5 (Leve Item ID 1681)
$ws.Range("H5").Value = 105
$ws.Range("J5").Value = 105
$ws.Range("L5").Value = 105
$ws.Range("N5").Value = -329
# Row 86 (Leve Item ID 11034)
$ws.Range("H86").Value = 14000
$ws.Range("J86").Value = 14000
$ws.Range("L86").Value = 14000
$ws.Range("N86").Value = -16372
# Row 89 (Leve Item ID 11034)
$ws.Range("H89").Value = 14000
$ws.Range("J89").Value = 14000
$ws.Range("L89").Value = 42000
$ws.Range("N89").Value = -53856

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 6969
$ws.Range("I40").Value = 7644.8887
$ws.Range("J40").Value = 6100
$ws.Range("K40").Value = 7644.8887
$ws.Range("L40").Value = 6100
$ws.Range("M40").Value = -7508.8887
$ws.Range("N40").Value = -6372
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 506.3846
$ws.Range("I55").Value = 297.875
$ws.Range("J55").Value = 840
$ws.Range("K55").Value = 297.875
$ws.Range("L55").Value = 840
$ws.Range("M55").Value = -124.875
$ws.Range("N55").Value = -1186
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 7358933.5
$ws.Range("I132").Value = 3408.9424
$ws.Range("J132").Value = 31264388
$ws.Range("K132").Value = 10226.8272
$ws.Range("L132").Value = 93793164
$ws.Range("M132").Value = -7696.8272
$ws.Range("N132").Value = -93798224
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 20005904
$ws.Range("I136").Value = 23811372
$ws.Range("J136").Value = 27201.25
$ws.Range("K136").Value = 71434116
$ws.Range("L136").Value = 81603.75
$ws.Range("M136").Value = -71431566
$ws.Range("N136").Value = -86703.75

$ws = $wb.Worksheets.Item("WVR")
# Row 14 (Leve Item ID 2658)
$ws.Range("H14").Value = 19431.5
$ws.Range("I14").Value = 23263
$ws.Range("J14").Value = 15600
$ws.Range("K14").Value = 23263
$ws.Range("L14").Value = 15600
$ws.Range("M14").Value = -23095
$ws.Range("N14").Value = -15936
